# Updated cryptos list values/links per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.224.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.769.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.98%  '
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("D4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5214'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +10.20%  '
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3611'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.25%  '
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07339'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.081'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.40%  '
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9987'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.059'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.10%  '
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.764.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.60%  '
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.959'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001042'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06406'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9980'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.842'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.268.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.09%  '
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.065'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.63%  '
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.322'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +11.00%  '
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.966.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.41%  '
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.058'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09757'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.46%  '
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.509'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.42%  '
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02218'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05973'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.829'
$ws.Range("D39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6102'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.430'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.980'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.67%  '
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.127'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5736'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.614'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.32%  '
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.877'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06693'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.37%  '
$ws.Range("E51").Style = "Normal"
